$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54, shifting existing rows 54-81 down to 55-82
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the latest weekly price record
$ws.Cells.Item(54, 1).Value = 8
$ws.Cells.Item(54, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44460
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(54, 6).Value = 100112044
$ws.Cells.Item(54, 7).Value = "Perejil"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 3000
$ws.Cells.Item(54, 11).Value = 1500
$ws.Cells.Item(54, 12).Value = 2000
$ws.Cells.Item(54, 13).Value = 1750
$ws.Cells.Item(54, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 16).Value = 1167
$ws.Cells.Item(54, 17).Value = 1.5
$ws.Cells.Item(54, 18).Value = "Hortaliza"
